$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCPAUrls")
$ws.Activate()
$win = $excel.ActiveWindow
Write-Host "win:" $win
Write-Host "FreezePanes:" $win.FreezePanes
Write-Host "SplitRow:" $win.SplitRow
